$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2307.3333
$ws.Range("J40").Value = 2666.3333
$ws.Range("L40").Value = 2666.3333
$ws.Range("N40").Value = -3016.3333

# Row 64
$ws.Range("H64").Value = 5085.7144
$ws.Range("J64").Value = 3866.6667
$ws.Range("L64").Value = 3866.6667
$ws.Range("N64").Value = -4362.6667

# Row 67
$ws.Range("H67").Value = 5085.7144
$ws.Range("J67").Value = 3866.6667
$ws.Range("L67").Value = 3866.6667
$ws.Range("N67").Value = -5582.6667

# Row 92
$ws.Range("H92").Value = 624.7083
$ws.Range("I92").Value = 527.2222
$ws.Range("K92").Value = 527.2222
$ws.Range("M92").Value = 720.7778

# Row 137
$ws.Range("H137").Value = 2975.625
$ws.Range("I137").Value = 1372.1333
$ws.Range("J137").Value = 4390.4707
$ws.Range("K137").Value = 4116.3999
$ws.Range("L137").Value = 13171.4121
$ws.Range("M137").Value = -1566.3999
$ws.Range("N137").Value = -18271.4121

# Row 138
$ws.Range("H138").Value = 3287.1096
$ws.Range("I138").Value = 1961.875
$ws.Range("J138").Value = 3450.2153
$ws.Range("K138").Value = 5885.625
$ws.Range("L138").Value = 10350.6459
$ws.Range("M138").Value = -745.625
$ws.Range("N138").Value = -20630.6459

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12479.585
$ws.Range("I32").Value = 9445.272000000001
$ws.Range("J32").Value = 24996.125
$ws.Range("K32").Value = 9445.272000000001
$ws.Range("L32").Value = 24996.125
$ws.Range("M32").Value = -9158.272000000001
$ws.Range("N32").Value = -25570.125

# Row 61
$ws.Range("H61").Value = 1972.75
$ws.Range("I61").Value = 1972.75
$ws.Range("K61").Value = 1972.75
$ws.Range("M61").Value = -1760.75

# Row 74
$ws.Range("H74").Value = 1748.6666
$ws.Range("I74").Value = 1158.3125
$ws.Range("K74").Value = 1158.3125
$ws.Range("M74").Value = -284.3125

# Row 77
$ws.Range("H77").Value = 1748.6666
$ws.Range("I77").Value = 1158.3125
$ws.Range("K77").Value = 5791.5625
$ws.Range("M77").Value = -1423.5625

# Row 122
$ws.Range("H122").Value = 2410.3333
$ws.Range("I122").Value = 1644.2941
$ws.Range("K122").Value = 4932.8823
$ws.Range("M122").Value = -2482.8823

# Row 132
$ws.Range("H132").Value = 2730.9092
$ws.Range("I132").Value = 2730.9092
$ws.Range("K132").Value = 8192.7276
$ws.Range("M132").Value = -5662.7276

# Row 136
$ws.Range("H136").Value = 1972.75
$ws.Range("I136").Value = 1972.75
$ws.Range("K136").Value = 5918.25
$ws.Range("M136").Value = -3368.25

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 298.15384
$ws.Range("J80").Value = 41.5
$ws.Range("L80").Value = 41.5
$ws.Range("N80").Value = -2037.5

# Row 83
$ws.Range("H83").Value = 298.15384
$ws.Range("J83").Value = 41.5
$ws.Range("L83").Value = 207.5
$ws.Range("N83").Value = -10191.5

# Row 105
$ws.Range("H105").Value = 2797.2856
$ws.Range("I105").Value = 2081.4688
$ws.Range("J105").Value = 10432.667
$ws.Range("K105").Value = 2081.4688
$ws.Range("L105").Value = 10432.667
$ws.Range("M105").Value = -334.4688000000001
$ws.Range("N105").Value = -13926.667

# Row 134
$ws.Range("H134").Value = 3960.5557
$ws.Range("I134").Value = 2938.6667
$ws.Range("K134").Value = 8816.000100000001
$ws.Range("M134").Value = -6281.000100000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4606.645
$ws.Range("I31").Value = 3039.7222
$ws.Range("K31").Value = 3039.7222
$ws.Range("M31").Value = -2744.7222

# Row 34
$ws.Range("H34").Value = 4606.645
$ws.Range("I34").Value = 3039.7222
$ws.Range("K34").Value = 3039.7222
$ws.Range("M34").Value = -2837.7222

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 513.7059
$ws.Range("J5").Value = 486.5
$ws.Range("L5").Value = 1459.5
$ws.Range("N5").Value = -1683.5

# Row 113
$ws.Range("H113").Value = 721.5
$ws.Range("I113").Value = 733.3333
$ws.Range("J113").Value = 718.7692
$ws.Range("K113").Value = 2199.9999
$ws.Range("L113").Value = 2156.3076
$ws.Range("M113").Value = -29.9998999999998
$ws.Range("N113").Value = -6496.3076

# Row 117
$ws.Range("H117").Value = 1402.2
$ws.Range("J117").Value = 3994
$ws.Range("L117").Value = 11982
$ws.Range("N117").Value = -18866

# Row 121
$ws.Range("H121").Value = 1193.75
$ws.Range("J121").Value = 1925
$ws.Range("L121").Value = 5775
$ws.Range("N121").Value = -8395

# Row 122
$ws.Range("H122").Value = 317.06668
$ws.Range("I122").Value = 279.4
$ws.Range("J122").Value = 392.4
$ws.Range("K122").Value = 2514.6
$ws.Range("L122").Value = 3531.6
$ws.Range("M122").Value = -64.59999999999991
$ws.Range("N122").Value = -8431.6

# Row 129
$ws.Range("H129").Value = 3498.375
$ws.Range("I129").Value = 2996.5
$ws.Range("J129").Value = 3665.6667
$ws.Range("K129").Value = 8989.5
$ws.Range("L129").Value = 10997.0001
$ws.Range("M129").Value = -3989.5
$ws.Range("N129").Value = -20997.0001

# Row 131
$ws.Range("H131").Value = 1508
$ws.Range("I131").Value = 1320.75
$ws.Range("J131").Value = 1561.5
$ws.Range("K131").Value = 3962.25
$ws.Range("L131").Value = 4684.5
$ws.Range("M131").Value = 1077.75
$ws.Range("N131").Value = -14764.5

# Row 135
$ws.Range("H135").Value = 513.7059
$ws.Range("J135").Value = 486.5
$ws.Range("L135").Value = 4378.5
$ws.Range("N135").Value = -9448.5

# Row 137
$ws.Range("H137").Value = 4880.737
$ws.Range("J137").Value = 5804.3335
$ws.Range("L137").Value = 17413.0005
$ws.Range("N137").Value = -27613.0005

$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 127
$ws.Range("I13").Value = 127
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 127
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 12
$ws.Range("N13").ClearContents()

# Row 29
$ws.Range("H29").Value = 6676833.5
$ws.Range("I29").Value = 13338833
$ws.Range("J29").Value = 14833.333
$ws.Range("K29").Value = 13338833
$ws.Range("L29").Value = 14833.333
$ws.Range("M29").Value = -13338543
$ws.Range("N29").Value = -15413.333

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3522.5
$ws.Range("I7").Value = 2366.6667
$ws.Range("K7").Value = 2366.6667
$ws.Range("M7").Value = -2254.6667

# Row 22
$ws.Range("H22").Value = 2347.5
$ws.Range("I22").Value = 990
$ws.Range("J22").Value = 2800
$ws.Range("K22").Value = 990
$ws.Range("L22").Value = 2800
$ws.Range("M22").Value = -695
$ws.Range("N22").Value = -3390

# Row 27
$ws.Range("H27").Value = 2347.5
$ws.Range("I27").Value = 990
$ws.Range("J27").Value = 2800
$ws.Range("K27").Value = 990
$ws.Range("L27").Value = 2800
$ws.Range("M27").Value = -883
$ws.Range("N27").Value = -3014

# Row 46
$ws.Range("H46").Value = 1040.6364
$ws.Range("I46").Value = 954.6667
$ws.Range("J46").Value = 1143.8
$ws.Range("K46").Value = 954.6667
$ws.Range("L46").Value = 1143.8
$ws.Range("M46").Value = -766.6667
$ws.Range("N46").Value = -1519.8

# Row 61
$ws.Range("H61").Value = 1728.4667
$ws.Range("I61").Value = 1768.7778
$ws.Range("K61").Value = 1768.7778
$ws.Range("M61").Value = -1566.7778

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

# Row 113
$ws.Range("H113").Value = 1728.4667
$ws.Range("I113").Value = 1768.7778
$ws.Range("K113").Value = 1768.7778
$ws.Range("M113").Value = 401.2221999999999

# Row 122
$ws.Range("H122").Value = 4333.636
$ws.Range("I122").Value = 4407.8887
$ws.Range("K122").Value = 13223.6661
$ws.Range("M122").Value = -10773.6661

# Row 126
$ws.Range("H126").Value = 3522.5
$ws.Range("I126").Value = 2366.6667
$ws.Range("K126").Value = 7100.000100000001
$ws.Range("M126").Value = -4630.000100000001

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 3333666.8
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

# Row 4
$ws.Range("H4").Value = 5000375
$ws.Range("I4").Value = 20000000
$ws.Range("K4").Value = 20000000
$ws.Range("M4").Value = -19999887

# Row 96
$ws.Range("H96").Value = 1495.2
$ws.Range("I96").Value = 1495.6666
$ws.Range("J96").Value = 1494.5
$ws.Range("K96").Value = 1495.6666
$ws.Range("L96").Value = 1494.5
$ws.Range("M96").Value = -122.6666
$ws.Range("N96").Value = -4240.5

# Row 113
$ws.Range("H113").Value = 673.375
$ws.Range("I113").Value = 932.6
$ws.Range("J113").Value = 241.33333
$ws.Range("K113").Value = 2797.8
$ws.Range("L113").Value = 723.99999
$ws.Range("M113").Value = -627.8000000000002
$ws.Range("N113").Value = -5063.99999

# Row 136
$ws.Range("H136").Value = 6244.25
$ws.Range("I136").Value = 5851.2856
$ws.Range("K136").Value = 17553.8568
$ws.Range("M136").Value = -15003.8568
